$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so
# numeric-looking strings (e.g. "1.00", "7.63") are not silently
# converted to numbers by Excel when the value is assigned.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.761.47"
$ws.Range("E2").Value = "  +3.00%  "

$ws.Range("D3").Value = "2.515.13"
$ws.Range("E3").Value = "  +1.04%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "595.70"
$ws.Range("E5").Value = "  +1.50%  "

$ws.Range("D6").Value = "177.61"
$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("D9").Value = "2.513.42"
$ws.Range("E9").Value = "  +0.92%  "

$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  +10.73%  "

$ws.Range("D12").Value = "0.342"
$ws.Range("E12").Value = "  +0.79%  "

$ws.Range("E13").Value = "  +1.52%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "25.99"
$ws.Range("E14").Value = "  +1.25%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.956.92"
$ws.Range("E15").Value = "  +0.48%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000177"
$ws.Range("E16").Value = "  +3.67%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "69.247.04"
$ws.Range("E17").Value = "  +2.40%  "

$ws.Range("D18").Value = "2.486.03"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("D19").Value = "7.63"
$ws.Range("E19").Value = "  +1.44%  "

$ws.Range("D20").Value = "363.69"
$ws.Range("E20").Value = "  +3.98%  "

$ws.Range("D21").Value = "11.06"

$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").Value = "70.58"
$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("D26").Value = "9.10"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Value = "1.68"
$ws.Range("E27").Value = "  -3.60%  "

$ws.Range("D28").Value = "2.643.86"
$ws.Range("E28").Value = "  +1.10%  "

$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("D30").Value = "513.03"
$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("D31").Value = "0.0₃0894"
$ws.Range("E31").Value = "  -0.63%  "

$ws.Range("D32").Value = "7.77"
$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("D33").Value = "1.24"
$ws.Range("E33").Value = "  -1.83%  "

$ws.Range("D34").Value = "1.78"
$ws.Range("E34").Value = "  +0.71%  "

$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.19%  "

$ws.Range("D36").Value = "161.94"
$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("E37").Value = "  -2.11%  "

$ws.Range("D38").Value = "18.84"
$ws.Range("E38").Value = "  +2.99%  "

$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").Value = "  -1.94%  "

$ws.Range("D42").Value = "1.73"
$ws.Range("E42").Value = "  -0.93%  "

$ws.Range("D43").Value = "4.81"
$ws.Range("E43").Value = "  -0.59%  "

$ws.Range("E44").Value = "  -2.34%  "

$ws.Range("D45").Value = "2.33"
$ws.Range("E45").Value = "  -3.07%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "38.83"
$ws.Range("E46").Value = "  -0.57%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "150.08"
$ws.Range("E47").Value = "  +3.79%  "

$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "3.58"
$ws.Range("E48").Value = "  +1.96%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "0.515"
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0738"
$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("D51").Value = "0.0₆0251"
$ws.Range("E51").Value = "  -0.94%  "

